$d = $word.ActiveDocument

# Locate the target paragraph (the one about SceneManager / Scenes).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.IndexOf("SceneManager") -ge 0) {
        $target = $p
        break
    }
}

$paraStart = $target.Range.Start
$fullText  = $target.Range.Text

$anchor   = "lumières"
$anchorAt = $fullText.IndexOf($anchor)
$splitPos = $paraStart + $anchorAt + $anchor.Length

$insertion = " (seulement pour la 3D)"

# 1) Insert the new text right after "lumières" first (this keeps it merged
#    into the existing run for the moment).
$insertPoint = $d.Range($splitPos, $splitPos)
$insertPoint.InsertAfter($insertion)

# 2) Toggle a character formatting property on exactly the newly inserted
#    span so the engine is forced to materialise it as its own run, then
#    flip the property back off so the visible formatting is unchanged.
$newSpan = $d.Range($splitPos, $splitPos + $insertion.Length)
$newSpan.Font.Bold = 1
$newSpan.Font.Bold = 0
